$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: 2021, Switzerland ---
$ws.Cells.Item(13, 1).Value = 2021
$ws.Cells.Item(13, 2).Value = "Switzerland"
$ws.Cells.Item(13, 3).Value = 8636560
$ws.Cells.Item(13, 4).Value = 34454
$ws.Cells.Item(13, 5).Value = 34330
$ws.Cells.Item(13, 6).Value = 125
$ws.Cells.Item(13, 7).Value = -1850
$ws.Cells.Item(13, 8).Value = 2080

$ws.Cells.Item(13, 9).Formula = '="0.4%"'
$ws.Cells.Item(13, 9).Copy()
$ws.Cells.Item(13, 9).PasteSpecial(-4163)

$ws.Cells.Item(13, 10).Formula = '="-5.4%"'
$ws.Cells.Item(13, 10).Copy()
$ws.Cells.Item(13, 10).PasteSpecial(-4163)

$ws.Cells.Item(13, 11).Formula = '="6.1%"'
$ws.Cells.Item(13, 11).Copy()
$ws.Cells.Item(13, 11).PasteSpecial(-4163)

$excel.CutCopyMode = $false

# --- Row 14: 2021, Sweden (only Year/Country/Population known so far) ---
$ws.Cells.Item(14, 1).Value = 2021
$ws.Cells.Item(14, 2).Value = "Sweden"
$ws.Cells.Item(14, 3).Value = 10379295
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 6).Style = "Normal"
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(14, 8).Style = "Normal"
$ws.Cells.Item(14, 9).Style = "Normal"
$ws.Cells.Item(14, 10).Style = "Normal"
$ws.Cells.Item(14, 11).Style = "Normal"

# --- Row 15: 2021, Spain ---
$ws.Cells.Item(15, 1).Value = 2021
$ws.Cells.Item(15, 2).Value = "Spain"
$ws.Cells.Item(15, 3).Value = 47353706
$ws.Cells.Item(15, 4).Value = 234826
$ws.Cells.Item(15, 5).Value = 221664
$ws.Cells.Item(15, 6).Value = 13163
$ws.Cells.Item(15, 7).Value = -292
$ws.Cells.Item(15, 8).Value = 26482

$ws.Cells.Item(15, 9).Formula = '="5.9%"'
$ws.Cells.Item(15, 9).Copy()
$ws.Cells.Item(15, 9).PasteSpecial(-4163)

$ws.Cells.Item(15, 10).Formula = '="-0.1%"'
$ws.Cells.Item(15, 10).Copy()
$ws.Cells.Item(15, 10).PasteSpecial(-4163)

$ws.Cells.Item(15, 11).Formula = '="11.9%"'
$ws.Cells.Item(15, 11).Copy()
$ws.Cells.Item(15, 11).PasteSpecial(-4163)

$excel.CutCopyMode = $false
